# Fruta / hortaliza, semanal
#
# A new weekly observation was inserted into the "Pomelo" price series.
# In the sheet this shows up as a brand-new row 91 (date 2021-12-09 /
# serial 44539) with the rest of the table (previously rows 91-186)
# pushed down by one row to 92-187, growing the used range from
# A1:T186 to A1:T187.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 91..186 down to 92..187, opening up a blank row 91.
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new weekly observation.
$ws.Range("A91").Value = 4
$ws.Range("B91").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C91").Value = "Los Lagos"
$ws.Range("D91").Value = 44539
$ws.Range("E91").Value = 10
$ws.Range("F91").Value = "Fruta"
$ws.Range("G91").Value = 100102
$ws.Range("H91").Value = "Cítricos"
$ws.Range("I91").Value = 100102006
$ws.Range("J91").Value = "Pomelo"
$ws.Range("K91").Value = "Start Ruby"
$ws.Range("L91").Value = "Primera"
$ws.Range("M91").Value = 200
$ws.Range("N91").Value = 11000
$ws.Range("O91").Value = 12000
$ws.Range("P91").Value = 11500
$ws.Range("Q91").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R91").Value = "Región de O'Higgins"
$ws.Range("S91").Value = 821
$ws.Range("T91").Value = 14
